$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns F and G
$ws.Range("F1").Value = "Generate Payments"
$ws.Range("G1").Value = "Payments Paid"

# Row 2
$ws.Range("F2").Value = "Yes"
$ws.Range("G2").Value = "Yes"

# Row 3
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = "No"

# Row 4
$ws.Range("F4").Value = "No"
$ws.Range("G4").Value = "No"

# Update column widths to match bestFit auto-widths (calibrated for the
# engine's pixel-rounding so the stored XML "width" lands as close as
# possible to the Excel-computed bestFit values of 16.875 / 12.875)
$ws.Columns("F").ColumnWidth = 16.142857142857142
$ws.Columns("G").ColumnWidth = 12.142857142857142

# Update selection to match diff target
$ws.Range("F5").Select()
